# 自动更新价格数据: insert today's row (2026-01-03) at the top of the
# data table, pushing all existing date rows down by one, which mirrors
# the upstream scraper re-running and prepending the newest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right below the header, shifting rows 2..44 down
# to rows 3..45 (dimension grows from A1:D44 to A1:D45).
$ws.Rows.Item(2).Insert()

# The inserted row inherits the header row's formatting (bold/border/
# centered). Strip that back to the plain, unstyled look the other data
# rows use.
$ws.Rows.Item(2).ClearFormats()

# Column A holds the date as literal text (e.g. "2026-01-03"), matching
# every other row in the sheet. Force text formatting first so Excel
# doesn't silently convert the string into a date serial number, then
# drop back to the workbook's default ("Normal") style once the literal
# text value is committed so no stray number format lingers on the cell.
$dateCell = $ws.Range("A2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-03"
$dateCell.Style = "Normal"

# The commodity price columns stay flat day over day.
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
